$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Paragraphs.Item(1).Range.Find.Execute("2025-07-12 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-13 Sunday", 2) | Out-Null

# Update the division problems/answers table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "62÷2=31, 0"
$t.Cell(1,2).Range.Text = "52÷2=26, 0"
$t.Cell(1,3).Range.Text = "71÷9=7, 8"
$t.Cell(1,4).Range.Text = "47÷9=5, 2"
$t.Cell(1,5).Range.Text = "90÷3=30, 0"

$t.Cell(5,1).Range.Text = "22÷5=4, 2"
$t.Cell(5,2).Range.Text = "95÷8=11, 7"
$t.Cell(5,3).Range.Text = "15÷5=3, 0"
$t.Cell(5,4).Range.Text = "27÷2=13, 1"
$t.Cell(5,5).Range.Text = "87÷7=12, 3"

$t.Cell(9,1).Range.Text = "95÷9=10, 5"
$t.Cell(9,2).Range.Text = "46÷2=23, 0"
$t.Cell(9,3).Range.Text = "14÷3=4, 2"
$t.Cell(9,4).Range.Text = "54÷7=7, 5"
$t.Cell(9,5).Range.Text = "36÷6=6, 0"

$t.Cell(13,1).Range.Text = "65÷2=32, 1"
$t.Cell(13,2).Range.Text = "93÷8=11, 5"
$t.Cell(13,3).Range.Text = "34÷4=8, 2"
$t.Cell(13,4).Range.Text = "31÷4=7, 3"
$t.Cell(13,5).Range.Text = "52÷8=6, 4"

$t.Cell(17,1).Range.Text = "25÷6=4, 1"
$t.Cell(17,2).Range.Text = "98÷9=10, 8"
$t.Cell(17,3).Range.Text = "29÷3=9, 2"
$t.Cell(17,4).Range.Text = "41÷4=10, 1"
$t.Cell(17,5).Range.Text = "83÷3=27, 2"

